$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.620.99'
$ws.Range("D3").Value = '3.026.02'
$ws.Range("E3").Value = '  -1.45%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '584.60'
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").Value = '147.43'
$ws.Range("E6").Value = '  -4.92%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.025.85'
$ws.Range("E8").Value = '  -1.39%  '
$ws.Range("E9").Value = '  -2.74%  '
$ws.Range("E10").Value = '  -3.84%  '
$ws.Range("D11").Value = '5.76'
$ws.Range("E11").Value = '  -1.45%  '
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("E13").Value = '  -2.57%  '
$ws.Range("D14").Value = '34.82'
$ws.Range("E14").Value = '  -5.50%  '
$ws.Range("E15").Value = '  +2.33%  '
$ws.Range("D16").Value = '3.525.66'
$ws.Range("E16").Value = '  -1.43%  '
$ws.Range("D17").Value = '7.06'
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("D18").Value = '62.582.43'
$ws.Range("E18").Value = '  -1.49%  '
$ws.Range("D19").Value = '3.025.76'
$ws.Range("E19").Value = '  -1.48%  '
$ws.Range("D20").Value = '464.84'
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("D21").Value = '13.98'
$ws.Range("E21").Value = '  -2.20%  '
$ws.Range("E22").Value = '  -2.16%  '
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("E24").Value = '  -3.19%  '
$ws.Range("D25").Value = '80.27'
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").Value = '12.44'
$ws.Range("E26").Value = '  -2.57%  '
$ws.Range("D27").Value = '10.07'
$ws.Range("E27").Value = '  -2.78%  '
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("D31").Value = '7.15'
$ws.Range("E31").Value = '  -3.61%  '
$ws.Range("D32").Value = '2.12'
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").Value = '27.50'
$ws.Range("E33").Value = '  +1.42%  '
$ws.Range("E34").Value = '  -3.73%  '
$ws.Range("E35").Value = '  -0.92%  '
$ws.Range("E36").Value = '  -2.90%  '
$ws.Range("E37").Value = '  -3.18%  '
$ws.Range("E38").Value = '  -3.09%  '
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("D40").Value = '9.04'
$ws.Range("E40").Value = '  -1.42%  '
$ws.Range("E41").Value = '  -10.15%  '
$ws.Range("D42").Value = '420.20'
$ws.Range("E42").Value = '  -3.63%  '
$ws.Range("E43").Value = '  +1.20%  '
$ws.Range("D45").Value = '2.787.65'
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("E46").Value = '  -1.14%  '
$ws.Range("D47").Value = '37.93'
$ws.Range("E47").Value = '  -6.05%  '
$ws.Range("D48").Value = '129.71'
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("E50").Value = '  -0.58%  '
$ws.Range("D51").Value = '24.12'
$ws.Range("E51").Value = '  -3.52%  '
